# Scheduled-runner market data refresh for the Rafflesia profit sheets.
# Updates currentAveragePrice / *PriceNQ / *PriceHQ / Leve profit columns
# (H:N) on the rows whose market snapshot changed, per sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 2777.3333
$ws.Range("J29").Value = 2777.3333
$ws.Range("L29").Value = 8331.999899999999
$ws.Range("N29").Value = -8893.999899999999

$ws.Range("H31").Value = 69
$ws.Range("I31").Value = 69
$ws.Range("K31").Value = 207
$ws.Range("M31").Value = 23

$ws.Range("H33").Value = 357.2
$ws.Range("I33").Value = 357.2
$ws.Range("K33").Value = 357.2
$ws.Range("M33").Value = -128.2

$ws.Range("H42").Value = 211.25
$ws.Range("I42").Value = 268.83334
$ws.Range("J42").Value = 153.66667
$ws.Range("K42").Value = 806.5000200000001
$ws.Range("L42").Value = 461.00001
$ws.Range("M42").Value = -576.5000200000001
$ws.Range("N42").Value = -921.00001

$ws.Range("H58").Value = 625.5
$ws.Range("J58").Value = 2232
$ws.Range("L58").Value = 6696
$ws.Range("N58").Value = -6996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1999
$ws.Range("I45").Value = 1999
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1999
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1622
$ws.Range("N45").ClearContents()

$ws.Range("H102").Value = 2462.5
$ws.Range("I102").Value = 2462.5
$ws.Range("K102").Value = 2462.5
$ws.Range("M102").Value = -840.5

$ws.Range("H103").Value = 30000
$ws.Range("J103").Value = 30000
$ws.Range("L103").Value = 30000
$ws.Range("N103").Value = -32344

$ws.Range("H110").Value = 810.4
$ws.Range("I110").Value = 763
$ws.Range("K110").Value = 763
$ws.Range("M110").Value = 1282

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 5806.2
$ws.Range("I132").Value = 5609.222
$ws.Range("K132").Value = 16827.666
$ws.Range("M132").Value = -14297.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1025.5
$ws.Range("I64").Value = 597.25
$ws.Range("J64").Value = 1311
$ws.Range("K64").Value = 597.25
$ws.Range("L64").Value = 1311
$ws.Range("M64").Value = -372.25
$ws.Range("N64").Value = -1761

$ws.Range("H67").Value = 1025.5
$ws.Range("I67").Value = 597.25
$ws.Range("J67").Value = 1311
$ws.Range("K67").Value = 597.25
$ws.Range("L67").Value = 1311
$ws.Range("M67").Value = 182.75
$ws.Range("N67").Value = -2871

$ws.Range("H94").Value = 3052.5
$ws.Range("I94").Value = 1400
$ws.Range("J94").Value = 3603.3333
$ws.Range("K94").Value = 1400
$ws.Range("L94").Value = 3603.3333
$ws.Range("M94").Value = -949
$ws.Range("N94").Value = -4505.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 250
$ws.Range("I22").Value = 250
$ws.Range("K22").Value = 250
$ws.Range("M22").Value = 100

$ws.Range("H31").Value = 5045.125
$ws.Range("I31").Value = 4999.6665
$ws.Range("J31").Value = 5072.4
$ws.Range("K31").Value = 4999.6665
$ws.Range("L31").Value = 5072.4
$ws.Range("M31").Value = -4704.6665
$ws.Range("N31").Value = -5662.4

$ws.Range("H34").Value = 5045.125
$ws.Range("I34").Value = 4999.6665
$ws.Range("J34").Value = 5072.4
$ws.Range("K34").Value = 4999.6665
$ws.Range("L34").Value = 5072.4
$ws.Range("M34").Value = -4797.6665
$ws.Range("N34").Value = -5476.4

$ws.Range("H86").Value = 5999.6665
$ws.Range("I86").Value = 4999.5
$ws.Range("J86").Value = 8000
$ws.Range("K86").Value = 4999.5
$ws.Range("L86").Value = 8000
$ws.Range("M86").Value = -3876.5
$ws.Range("N86").Value = -10246

$ws.Range("H89").Value = 5999.6665
$ws.Range("I89").Value = 4999.5
$ws.Range("J89").Value = 8000
$ws.Range("K89").Value = 24997.5
$ws.Range("L89").Value = 40000
$ws.Range("M89").Value = -19381.5
$ws.Range("N89").Value = -51232

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()

$ws.Range("H12").Value = 29.25
$ws.Range("J12").Value = 33.5
$ws.Range("L12").Value = 100.5
$ws.Range("N12").Value = -446.5

$ws.Range("H14").Value = 2498.5
$ws.Range("I14").Value = 2498.5
$ws.Range("K14").Value = 7495.5
$ws.Range("M14").Value = -7322.5

$ws.Range("H26").Value = 30
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

$ws.Range("H33").Value = 103
$ws.Range("I33").Value = 103.6
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 621.5999999999999
$ws.Range("L33").Value = 600
$ws.Range("M33").Value = -338.5999999999999
$ws.Range("N33").Value = -1166

$ws.Range("H34").Value = 2987.5
$ws.Range("J34").Value = 2987.5
$ws.Range("L34").Value = 8962.5
$ws.Range("N34").Value = -9130.5

$ws.Range("H39").Value = 1552.7778
$ws.Range("J39").Value = 2993.75
$ws.Range("L39").Value = 8981.25
$ws.Range("N39").Value = -9569.25

$ws.Range("H55").Value = 2987.5
$ws.Range("J55").Value = 2987.5
$ws.Range("L55").Value = 8962.5
$ws.Range("N55").Value = -9316.5

$ws.Range("H121").Value = 844.5
$ws.Range("J121").Value = 844.5
$ws.Range("L121").Value = 2533.5
$ws.Range("N121").Value = -5153.5

$ws.Range("H131").Value = 2601.5454
$ws.Range("J131").Value = 2887.4443
$ws.Range("L131").Value = 8662.332900000001
$ws.Range("N131").Value = -18742.3329

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1566.6666
$ws.Range("I122").Value = 1566.6666
$ws.Range("K122").Value = 4699.9998
$ws.Range("M122").Value = -2249.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 999
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 999
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 999
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1589

$ws.Range("H27").Value = 999
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 999
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 999
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -1213

$ws.Range("H122").Value = 15067.25
$ws.Range("I122").Value = 15067.25
$ws.Range("K122").Value = 45201.75
$ws.Range("M122").Value = -42751.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 606.1667
$ws.Range("I122").Value = 606.1667
$ws.Range("K122").Value = 1818.5001
$ws.Range("M122").Value = 631.4999

$ws.Range("H124").Value = 94333.336
$ws.Range("J124").Value = 94333.336
$ws.Range("L124").Value = 94333.336
$ws.Range("N124").Value = -104153.336
